$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(8, 1).Value = "timer_activate"
$ws.Cells.Item(8, 2).Value = "TIMER"
$ws.Cells.Item(8, 3).NumberFormat = "@"
$ws.Cells.Item(8, 3).Value = "32"
$ws.Cells.Item(8, 3).ClearFormats()
$ws.Cells.Item(8, 4).Value = 1
$ws.Cells.Item(8, 5).Value = "6/6"

$ws.Cells.Item(9, 1).Value = "semaphore_ceiling_put"
$ws.Cells.Item(9, 2).Value = "SEMAPHORE"
$ws.Cells.Item(9, 3).Value = "8,12,16"
$ws.Cells.Item(9, 4).Value = 3
$ws.Cells.Item(9, 5).Value = "3/3"

$ws.Cells.Item(10, 1).Value = "block_allocate"
$ws.Cells.Item(10, 2).Value = "BLOCK"
$ws.Cells.Item(10, 3).Value = "8,16,36"
$ws.Cells.Item(10, 4).Value = 3
$ws.Cells.Item(10, 5).Value = "2/2"

$ws.Cells.Item(11, 1).Value = "timer_change"
$ws.Cells.Item(11, 2).Value = "TIMER"
$ws.Cells.Item(11, 3).Value = "8,12"
$ws.Cells.Item(11, 4).Value = 2
$ws.Cells.Item(11, 5).Value = "2/2"

$ws.Cells.Item(12, 1).Value = "mutex_put"
$ws.Cells.Item(12, 2).Value = "MUTEX"
$ws.Cells.Item(12, 3).Value = "8,20,24,28,40"
$ws.Cells.Item(12, 4).Value = 5
$ws.Cells.Item(12, 5).Value = "8/3"

$ws.Cells.Item(13, 1).Value = "event_flags_set_notify"
$ws.Cells.Item(13, 2).Value = "EVENT"
$ws.Cells.Item(13, 3).Value = "40,44"
$ws.Cells.Item(13, 4).Value = 2
$ws.Cells.Item(13, 5).Value = "1/1"

$ws.Cells.Item(14, 1).Value = "thread_priority_change"
$ws.Cells.Item(14, 2).Value = "THREAD"
$ws.Cells.Item(14, 3).Value = "32,36,44,60,204,208"
$ws.Cells.Item(14, 4).Value = 6
$ws.Cells.Item(14, 5).Value = "7/3"

$ws.Cells.Item(15, 1).Value = "semaphore_put"
$ws.Cells.Item(15, 2).Value = "SEMAPHORE"
$ws.Cells.Item(15, 3).Value = "12,16"
$ws.Cells.Item(15, 4).Value = 2
$ws.Cells.Item(15, 5).Value = "3/3"

$ws.Cells.Item(16, 1).Value = "timer_deactivate"
$ws.Cells.Item(16, 2).Value = "TIMER"
$ws.Cells.Item(16, 3).NumberFormat = "@"
$ws.Cells.Item(16, 3).Value = "8"
$ws.Cells.Item(16, 3).ClearFormats()
$ws.Cells.Item(16, 4).Value = 1
$ws.Cells.Item(16, 5).Value = "6/4"

$ws.Cells.Item(17, 1).Value = "queue_flush"
$ws.Cells.Item(17, 2).Value = "QUEUE"
$ws.Cells.Item(17, 3).Value = "20,32,36"
$ws.Cells.Item(17, 4).Value = 3
$ws.Cells.Item(17, 5).Value = "2/2"

$ws.Cells.Item(18, 1).Value = "event_flags_create"
$ws.Cells.Item(18, 2).Value = "EVENT"
$ws.Cells.Item(18, 3).NumberFormat = "@"
$ws.Cells.Item(18, 3).Value = "4"
$ws.Cells.Item(18, 3).ClearFormats()
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(18, 5).Value = "1/1"

$ws.Cells.Item(19, 1).Value = "mutex_delete"
$ws.Cells.Item(19, 2).Value = "MUTEX"
$ws.Cells.Item(19, 3).Value = "20,24,40"
$ws.Cells.Item(19, 4).Value = 3
$ws.Cells.Item(19, 5).Value = "4/3"

$ws.Cells.Item(20, 1).Value = "event_flags_set"
$ws.Cells.Item(20, 2).Value = "EVENT"
$ws.Cells.Item(20, 3).Value = "8,20,32"
$ws.Cells.Item(20, 4).Value = 3
$ws.Cells.Item(20, 5).Value = "27/1"

$ws.Cells.Item(21, 1).Value = "thread_create"
$ws.Cells.Item(21, 2).Value = "THREAD"
$ws.Cells.Item(21, 3).Value = "16,20,24,28,40,44,60,68,72,152,156,184,188,204,208"
$ws.Cells.Item(21, 4).Value = 15
$ws.Cells.Item(21, 5).Value = "1/1"

$ws.Cells.Item(22, 1).Value = "block_pool_create"
$ws.Cells.Item(22, 2).Value = "BLOCK"
$ws.Cells.Item(22, 3).Value = "4,20,24,28"
$ws.Cells.Item(22, 4).Value = 4
$ws.Cells.Item(22, 5).Value = "1/1"

$ws.Cells.Item(23, 1).Value = "event_flags_get"
$ws.Cells.Item(23, 2).Value = "EVENT"
$ws.Cells.Item(23, 3).Value = "8,20,32"
$ws.Cells.Item(23, 4).Value = 3
$ws.Cells.Item(23, 5).Value = "6/4"

$ws.Cells.Item(24, 1).Value = "semaphore_create"
$ws.Cells.Item(24, 2).Value = "SEMAPHORE"
$ws.Cells.Item(24, 3).Value = "4,8"
$ws.Cells.Item(24, 4).Value = 2
$ws.Cells.Item(24, 5).Value = "1/1"

$ws.Cells.Item(25, 1).Value = "thread_reset"
$ws.Cells.Item(25, 2).Value = "THREAD"
$ws.Cells.Item(25, 3).NumberFormat = "@"
$ws.Cells.Item(25, 3).Value = "8"
$ws.Cells.Item(25, 3).ClearFormats()
$ws.Cells.Item(25, 4).Value = 1
$ws.Cells.Item(25, 5).Value = "2/2"

$ws.Cells.Item(26, 1).Value = "semaphore_put_notify"
$ws.Cells.Item(26, 2).Value = "SEMAPHORE"
$ws.Cells.Item(26, 3).Value = "32,36"
$ws.Cells.Item(26, 4).Value = 2
$ws.Cells.Item(26, 5).Value = "1/1"

$ws.Cells.Item(27, 1).Value = "timer_create"
$ws.Cells.Item(27, 2).Value = "TIMER"
$ws.Cells.Item(27, 3).Value = "4,8,12,20,32,44,48"
$ws.Cells.Item(27, 4).Value = 7
$ws.Cells.Item(27, 5).Value = "2/1"

$ws.Cells.Item(28, 1).Value = "byte_pool_create"
$ws.Cells.Item(28, 2).Value = "BYTE"
$ws.Cells.Item(28, 3).Value = "4,8,16,20,24,28"
$ws.Cells.Item(28, 4).Value = 6
$ws.Cells.Item(28, 5).Value = "1/1"

$ws.Cells.Item(29, 1).Value = "queue_send_notify"
$ws.Cells.Item(29, 2).Value = "QUEUE"
$ws.Cells.Item(29, 3).Value = "60,64"
$ws.Cells.Item(29, 4).Value = 2
$ws.Cells.Item(29, 5).Value = "1/1"

$ws.Cells.Item(30, 1).Value = "queue_front_send"
$ws.Cells.Item(30, 2).Value = "QUEUE"
$ws.Cells.Item(30, 3).Value = "32,40,44"
$ws.Cells.Item(30, 4).Value = 3
$ws.Cells.Item(30, 5).Value = "4/3"

$ws.Cells.Item(31, 1).Value = "queue_receive"
$ws.Cells.Item(31, 2).Value = "QUEUE"
$ws.Cells.Item(31, 3).Value = "16,32,36,40,44"
$ws.Cells.Item(31, 4).Value = 5
$ws.Cells.Item(31, 5).Value = "7/3"

$ws.Cells.Item(32, 1).Value = "thread_preemption_change"
$ws.Cells.Item(32, 2).Value = "THREAD"
$ws.Cells.Item(32, 3).NumberFormat = "@"
$ws.Cells.Item(32, 3).Value = "60,208"
$ws.Cells.Item(32, 3).ClearFormats()
$ws.Cells.Item(32, 4).Value = 2
$ws.Cells.Item(32, 5).Value = "5/3"
